$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Chap" column (A) values for rows 2-4 from 1 to 2
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 2

# Update the view/selection state: scroll to A4 and select G4
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G4").Select()
